$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# QuestType data (column C) added for rows 2-30
$questTypes = @{
    2 = 2
    3 = 3
    4 = 0
    5 = 0
    6 = 0
    7 = 5
    8 = 5
    9 = 5
    10 = 5
    11 = 5
    12 = 5
    13 = 5
    14 = 5
    15 = 0
    16 = 1
    17 = 5
    18 = 3
    19 = 1
    20 = 3
    21 = 0
    22 = 5
    23 = 5
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 5
    29 = 5
    30 = 5
}

foreach ($row in $questTypes.Keys) {
    $ws.Cells.Item($row, 3).Value = $questTypes[$row]
}

# Unhide column C (QuestType) now that it has data
$ws.Columns("C").Hidden = $false

# Restore the view state captured at save time (best-effort; some
# window chrome like pixel size / GUIDs is host-generated and not
# settable from the object model).
$ws.Range("B30").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 2

